$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 = "Save", matching style of other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add new data cell H2 = 1 (numeric)
$ws.Range("H2").Value = 1
